# Update the public EPEX Spot prices workbook with the latest daily data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column AS (28-jul) with hourly prices.
# ---------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Header cell: copy the style from the previous day's header (AR1) so the
# new column keeps the same bold / centered / bordered look, then set its
# text.
$wsSpot.Range("AR1").Copy($wsSpot.Range("AS1"))
$wsSpot.Range("AS1").Value = "28-jul"

$spotValues = @{
    2  = 62.29
    3  = 55.75
    4  = 48.95
    5  = 19.96
    6  = 38
    7  = 45.07
    8  = 50
    9  = 56.23
    10 = 63.71
    11 = 30
    12 = 20.59
    13 = 22.49
    14 = 21.87
    15 = 20.22
    16 = 11.17
    17 = 10.87
    18 = 14.35
    19 = 22.97
    20 = 44.3
    21 = 60.11
    22 = 60
    23 = 74.98
    24 = 88.54000000000001
    25 = 63.25
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Range("AS$row").Value = $spotValues[$row]
}

# ---------------------------------------------------------------------
# Helper: append a (date, price) row as plain text / number, without
# letting Excel auto-convert the ISO date string into a date serial
# number (which would add an unwanted number format / style).
# ---------------------------------------------------------------------
function Add-DateRow {
    param($ws, [int]$row, [string]$dateText, [double]$price)

    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $dateText
    $ws.Range("A$row").ClearFormats()

    $ws.Range("B$row").Value = $price
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append rows 42 & 43.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
Add-DateRow $wsGaz 42 "2025-07-26" 31.85
Add-DateRow $wsGaz 43 "2025-07-27" 31.85

# ---------------------------------------------------------------------
# Sheet "CO2": append rows 42 & 43.
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
Add-DateRow $wsCo2 42 "2025-07-26" 70.7
Add-DateRow $wsCo2 43 "2025-07-27" 70.7
